$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update descriptive text labels (shared strings) ---
$ws.Range("B1").Value = "2014-10 Prueba 2"
$ws.Range("B18").Value = "Burger Xpress"
$ws.Range("B89").Value = "Método Gráfico"
$ws.Range("B326").Value = "Selección de Proyectos"
$ws.Range("B429").Value = "Bono 1"

# --- Update grade (column H) values for the affected rows ---
$hValues = @{
    23 = 0.5
    28 = 0.5
    32 = 0.4
    33 = 0.4
    38 = 0.4
    42 = 0.4
    51 = 1.8
    56 = 1.8
    63 = 1.5
    67 = 0.5
    72 = 1.5
    78 = 2.5
    87 = 1.5
    94 = 0.5
    95 = 0.5
    96 = 0.2
    101 = 0.5
    103 = 0.5
    104 = 0.5
    105 = 0.2
    107 = 0.5
    108 = 0.5
    110 = 0.5
    111 = 0.5
    116 = 0.5
    126 = 0.5
    127 = 0.5
    129 = 0.5
    130 = 0.5
    136 = 0.3
    141 = 0.5
    142 = 0.3
    144 = -0.5
    154 = 0.5
    157 = 0.5
    158 = 0.5
    160 = 0.5
    161 = 0.5
    163 = 0.5
    164 = 0.5
    167 = 0.4
    168 = 0.4
    169 = 0.2
    172 = 0.1
    173 = 0.1
    177 = 0.5
    178 = 0.5
    182 = 0.3
    183 = 0.3
    185 = 0.2
    186 = 0.2
    188 = 0.5
    189 = 0.5
    192 = 0.5
    193 = 0.5
    194 = 0.3
    199 = 0.5
    200 = 0.5
    204 = 0.3
    205 = 0.3
    207 = 0.2
    208 = 0.2
    210 = 0.5
    211 = 0.5
    214 = 0.5
    215 = 0.5
    216 = 0.3
    221 = 0.5
    222 = 0.5
    226 = 0.3
    227 = 0.3
    229 = 0.2
    230 = 0.2
    232 = 0.5
    233 = 0.5
    236 = 0.5
    237 = 0.5
    238 = 0.3
    250 = 0.5
    251 = 0.5
    256 = 1.5
    274 = 0.5
    275 = 0.5
    276 = 0.3
    278 = 0.5
    279 = 0.5
    280 = 0.3
    285 = 0.5
    303 = 0.5
    317 = 0.5
    321 = 0.5
    329 = 0.3
    330 = 0.3
    331 = 0.1
    334 = 0.3
    335 = 0.3
    336 = 0.1
    338 = 0.4
    339 = 0.4
    340 = 0.2
    341 = 0.1
    344 = 0.5
    345 = 0.5
    346 = 0.2
    348 = 0.5
    349 = 0.5
    350 = 0.5
    351 = 0.25
    356 = 0.5
    357 = 0.5
    361 = 0.5
    362 = 0.5
    366 = 0.5
    371 = 0.5
    375 = 1.5
    376 = 1.5
    383 = 1.5
    384 = 1.5
    385 = 0.7
    393 = 0.5
    396 = 0.5
    397 = 0.5
    399 = 0.5
    400 = 0.5
    404 = 0.5
    408 = 0.5
    413 = 0.5
    417 = 1.5
    418 = 1.5
    425 = 1.5
    426 = 1.5
}

foreach ($row in $hValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $hValues[$row]
}

Write-Host "Updated $($hValues.Count) grade cells and 5 labels"
